$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F, shifting the old District column (F) to G
$ws.Columns.Item(6).Insert()

# Populate the new "Address" column (F) with the address/taluk portion
# extracted from the combined name+address text already in columns B/E
$ws.Range("F2").Value = 'Address'
$ws.Range("F3").Value = 'G H S SonanaAurad'
$ws.Range("F4").Value = 'G(Ex-ML) P U College'
$ws.Range("F5").Value = 'Govt Girls High School WadgaonAurad'
$ws.Range("F7").Value = 'G H S DashavaraChannapatna'
$ws.Range("F8").Value = 'G H S MokaliArakalgud'
$ws.Range("F9").Value = 'Kondaji Basappa H S Betur'
$ws.Range("F10").Value = 'G G H S Channageri'
$ws.Range("F11").Value = 'Bhadravathi'
$ws.Range("F12").Value = 'S G H S GiriyapurKadur'
$ws.Range("F13").Value = 'Vijaya High School Santhebennurchannagiri'
$ws.Range("F14").Value = 'M M G H S M CampBirur'
$ws.Range("F15").Value = 'G H S HonganurChannapatna'
$ws.Range("F16").Value = 'Anand High School SanthpurAurad(B)'
$ws.Range("F18").Value = 'College RabakaviJamkhandi'
$ws.Range("F19").Value = 'Sree Jagathguru Renukacharya High School Uddeboranahalli'
$ws.Range("F20").Value = 'G B J C Arkalgud'
$ws.Range("F21").Value = 'S V S H S RudrapatnaArakalgud'
$ws.Range("F22").Value = 'K E B H S Malmadeli'
$ws.Range("F23").Value = 'Al Ameen Alhadad high SchoolChannapatna'
$ws.Range("F24").Value = 'S J C Channapatna'
$ws.Range("F25").Value = 'G H S HulikalArakalgud'
$ws.Range("F26").Value = 'S S High School NivagundaBadami'
$ws.Range("F27").Value = 'S P R H S RamanathapuraArakalgud'
$ws.Range("F28").Value = 'S A P U C (H S ) HalasangiChadchan'
$ws.Range("F29").Value = 'G H S Thalihalla'
$ws.Range("F30").Value = 'S S M S Chadchan'
$ws.Range("F31").Value = 'G H S AnkalagaAfzalpur'
$ws.Range("F32").Value = 'G J C DoddamaggeArakalgud'
$ws.Range("F34").Value = 'Hosanagar'
$ws.Range("F35").Value = 'G P U C Hosanagara'
$ws.Range("F36").Value = 'Bhadra High SchoolBhadravathi'
